$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on numeric-looking price cells so they remain stored as text
# (looping individually, since comma-separated Range union only applies to the first area)
$numericTextCells = @("D5", "D6", "D9", "D10", "D12", "D13", "D17", "D19", "D20", "D21", "D23", "D24", "D27", "D28", "D29", "D33", "D35", "D36", "D38", "D40", "D43", "D44", "D46", "D48", "D50", "D51")
foreach ($cellRef in $numericTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range('D2').Value = '51.848.28'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '2.810.39'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '355.96'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').Value = '111.76'
$ws.Range('E6').Value = '  +2.41%  '
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.630'
$ws.Range('E9').Value = '  +8.01%  '
$ws.Range('D10').Value = '40.36'
$ws.Range('E10').Value = '  +2.19%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '0.0841'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = '20.01'
$ws.Range('E13').Value = '  +2.97%  '
$ws.Range('E14').Value = '  +2.88%  '
$ws.Range('D15').Value = '3.252.69'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').Value = '2.810.26'
$ws.Range('E16').Value = '  +1.83%  '
$ws.Range('D17').Value = '0.942'
$ws.Range('E17').Value = '  +1.50%  '
$ws.Range('D18').Value = '51.831.65'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '7.66'
$ws.Range('E19').Value = '  +2.98%  '
$ws.Range('D20').Value = '3.21'
$ws.Range('E20').Value = '  +4.21%  '
$ws.Range('D21').Value = '13.65'
$ws.Range('E21').Value = '  +3.96%  '
$ws.Range('E22').Value = '  +1.36%  '
$ws.Range('D23').Value = '70.47'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').Value = '268.96'
$ws.Range('E24').Value = '  +0.61%  '
$ws.Range('E25').Value = '  +1.83%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '26.19'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').Value = '0.162'
$ws.Range('E28').Value = '  -1.13%  '
$ws.Range('D29').Value = '38.73'
$ws.Range('E29').Value = '  +12.70%  '
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('E32').Value = '  +2.01%  '
$ws.Range('D33').Value = '6.12'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  +9.53%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = '0.0445'
$ws.Range('E35').Value = '  +0.40%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '0.0882'
$ws.Range('E36').Value = '  +5.60%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = '18.84'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').Value = '3.14'
$ws.Range('E40').Value = '  +1.24%  '
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('E42').Value = '  -1.10%  '
$ws.Range('D43').Value = '120.72'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('D44').Value = '22.17'
$ws.Range('E44').Value = '  +1.89%  '
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('D46').Value = '3.42'
$ws.Range('E46').Value = '  +5.00%  '
$ws.Range('D47').Value = '2.112.19'
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +6.41%  '
$ws.Range('E49').Value = '  +1.25%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '1.37'
$ws.Range('E50').Value = '  +6.77%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '5.46'
$ws.Range('E51').Value = '  -1.38%  '
